# Update "Línea 141" horarios workbook with the latest scrape snapshot.
# New scrape timestamp: 02:16:52 (was 01:53:21)

$wb = $excel.ActiveWorkbook

$newTime = "02:16:52"

# --- Sheet 1: LP1912 (has the actual arrival-time data rows) ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: " + $newTime

$ws1.Range("A6").Value = $newTime
$ws1.Range("B6").Value = "03:02"
$ws1.Range("C6").Value = "15_ABASTO"
$ws1.Range("D6").Value = 46
$ws1.Range("E6").Value = "LP1912"

$ws1.Range("A7").Value = $newTime
$ws1.Range("B7").Value = "03:48"
$ws1.Range("C7").Value = "14_ABASTO"
$ws1.Range("D7").Value = 92
$ws1.Range("E7").Value = "LP1912"

$ws1.Range("A8").Value = $newTime
$ws1.Range("B8").Value = "04:01"
$ws1.Range("C8").Value = "81_EL PELIGRO"
$ws1.Range("D8").Value = 105
$ws1.Range("E8").Value = "LP1912"

# --- Sheet 2: LP1912-215 (no data rows, just refresh the timestamp) ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: " + $newTime

# --- Sheet 3: 6203-6173 (no data rows, just refresh the timestamp) ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: " + $newTime
